# Remove the "Requisitos" Heading2 paragraph and the requirement list
# paragraph that immediately follows it
# ("LOM3238 -  Projeto Integrado  (Requisito)").

$d = $word.ActiveDocument

$count = $d.Paragraphs.Count
$startPara = $null

for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $txt = $p.Range.Text
    # The "Requisitos" paragraph is a short heading consisting of just that
    # word (plus its trailing paragraph mark) - match on trimmed text so we
    # don't accidentally match a longer paragraph that merely starts with it.
    $trimmed = $txt.Trim()
    if ($trimmed -eq "Requisitos") {
        $startPara = $i
        break
    }
}

if ($null -ne $startPara -and ($startPara + 1) -le $count) {
    $first = $d.Paragraphs.Item($startPara)
    $last = $d.Paragraphs.Item($startPara + 1)

    $start = $first.Range.Start
    $end = $last.Range.End

    $r = $d.Range($start, $end)
    $r.Delete()
}
